$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.967.64"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.449.35"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'523.24"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'131.11"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "2.452.31"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "'0.0981"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "2.884.65"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "57.879.55"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'21.72"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "2.448.61"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").Value = "'4.12"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "'314.05"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'64.96"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").Value = "'173.17"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").Value = "'6.22"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  -4.62%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D36").Value = "'17.78"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").Value = "'3.80"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "'36.23"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "'0.800"
$ws.Range("E41").Value = "  +4.29%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.586"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'261.01"
$ws.Range("E44").Value = "  -5.28%  "
$ws.Range("E45").Value = "  -5.57%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0920"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'122.33"
$ws.Range("E47").Value = "  -5.40%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'16.97"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("E51").Value = "  -3.56%  "
